{"js": "// Recolor the \"16 / NodeJS / Interacting With NodeJS\" syllabus row\n// (Key Features Of AngularJS) from black to green (00B050), matching\n// the color already used by the surrounding table rows.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Locate the row whose first cell reads \"16\" and whose second cell\n// reads \"NodeJS\" -- this is more robust than a fixed row index.\nlet targetRow = null;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  if (cells.items.length >= 2) {\n    cells.items[0].load(\"value\");\n    cells.items[1].load(\"value\");\n    await context.sync();\n    if (cells.items[0].value.trim() === \"16\" && cells.items[1].value.trim() === \"NodeJS\") {\n      targetRow = row;\n      break;\n    }\n  }\n}\n\nif (targetRow) {\n  const cells = targetRow.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    // Set the color on the whole cell body (paragraph mark + run text)\n    // so both the paragraph's rPr and the run's rPr pick up the new color,\n    // mirroring what Word does when you select an entire cell and apply\n    // a font color.\n    cell.body.font.color = \"#00B050\";\n  }\n  await context.sync();\n}\n", "ps1": "# Recolor the \"16 / NodeJS / Interacting With NodeJS\" syllabus row\n# (Key Features Of AngularJS) from black to green (00B050), matching\n# the color already used by the surrounding table rows.\n$doc = $word.ActiveDocument\n$tbl = $doc.Tables.Item(1)\n\n# Compute the Word (BGR) color value for hex RGB 00B050.\n$r = [Convert]::ToInt32(\"00\", 16)\n$g = [Convert]::ToInt32(\"B0\", 16)\n$b = [Convert]::ToInt32(\"50\", 16)\n$greenColor = $b * 65536 + $g * 256 + $r\n\n# Locate the row whose first cell reads \"16\" and whose second cell\n# reads \"NodeJS\" -- more robust than relying on a fixed row index.\n$targetRow = $null\nfor ($i = 1; $i -le $tbl.Rows.Count; $i++) {\n    $row = $tbl.Rows.Item($i)\n    if ($row.Cells.Count -ge 2) {\n        $cell1Text = $row.Cells.Item(1).Range.Text.Trim([char]7, [char]13)\n        $cell2Text = $row.Cells.Item(2).Range.Text.Trim([char]7, [char]13)\n        if ($cell1Text -eq \"16\" -and $cell2Text -eq \"NodeJS\") {\n            $targetRow = $row\n            break\n        }\n    }\n}\n\nif ($targetRow -ne $null) {\n    for ($i = 1; $i -le $targetRow.Cells.Count; $i++) {\n        $cell = $targetRow.Cells.Item($i)\n        # Set color on the whole cell range so both the paragraph mark's\n        # rPr and the run's rPr pick up the new color.\n        $cell.Range.Font.Color = $greenColor\n    }\n}\n"}
